$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Wireframes and Mockups for Unified Operations Dashboard for Managers"
$ws.Range("B1").Value = "As a Product Team, we need wireframes and mockups for the 'Unified Operations Dashboard for Managers' so that managers have a cohesive, intuitive, and user-centric experience from day one."
$ws.Range("C1").Value = "GIVEN a manager logs in, WHEN they view their dashboard, THEN all core modules are presented clearly and intuitively.`nGIVEN a user views any screen, WHEN they interact with it, THEN the design is consistent with the established brand and style guide."
$ws.Range("D1").Value = "Not specified"
$ws.Range("E1").Value = "Not specified"
$ws.Range("F1").Value = "None specified"
$ws.Range("G1").Value = "Scope limited to wireframes and mockups for the Unified Operations Dashboard for Managers."

$ws.Range("A2").Value = "Design for Integrated Incident Reporting Form and Submission Flow"
$ws.Range("B2").Value = "As a Product Team, we need the design for the 'Integrated Incident Reporting' form and submission flow so that staff can report incidents simply, quickly, and with minimal training."
$ws.Range("C2").Value = "GIVEN a staff member needs to report an incident, WHEN they access the form, THEN the process is simple, quick, and requires minimal training.`nGIVEN a user views any screen, WHEN they interact with it, THEN the design is consistent with the established brand and style guide."
$ws.Range("D2").Value = "Not specified"
$ws.Range("E2").Value = "Not specified"
$ws.Range("F2").Value = "None specified"
$ws.Range("G2").Value = "Scope limited to the design and submission flow for the Integrated Incident Reporting form."

$ws.Range("A3").Value = "UI for Staff Task & Schedule Viewer"
$ws.Range("B3").Value = "As a Product Team, we need the UI for the 'Staff Task & Schedule Viewer' so that staff can intuitively view their tasks and schedules."
$ws.Range("C3").Value = "GIVEN a user views any screen, WHEN they interact with it, THEN the design is consistent with the established brand and style guide."
$ws.Range("D3").Value = "Not specified"
$ws.Range("E3").Value = "Not specified"
$ws.Range("F3").Value = "None specified"
$ws.Range("G3").Value = "Scope limited to the UI for the Staff Task & Schedule Viewer."

$ws.Range("A4").Value = "Visual Design for Compliance & Certification Tracking Module"
$ws.Range("B4").Value = "As a Product Team, we need the visual design for the 'Compliance & Certification Tracking' module so that users can track compliance and certifications within a cohesive and branded interface."
$ws.Range("C4").Value = "GIVEN a user views any screen, WHEN they interact with it, THEN the design is consistent with the established brand and style guide."
$ws.Range("D4").Value = "Not specified"
$ws.Range("E4").Value = "Not specified"
$ws.Range("F4").Value = "None specified"
$ws.Range("G4").Value = "Scope limited to the visual design for the Compliance & Certification Tracking module."

$ws.Range("A5").Value = "User Flow Mapping for Closed-Loop Incident-to-Training Workflow"
$ws.Range("B5").Value = "As a Product Team, we need user flow mapping for the 'Closed-Loop Incident-to-Training Workflow' so that the process from incident reporting to training is clearly defined and user-centric."
$ws.Range("C5").Value = "GIVEN a user views any screen, WHEN they interact with it, THEN the design is consistent with the established brand and style guide."
$ws.Range("D5").Value = "Not specified"
$ws.Range("E5").Value = "Not specified"
$ws.Range("F5").Value = "None specified"
$ws.Range("G5").Value = "Scope limited to user flow mapping for the Closed-Loop Incident-to-Training Workflow."

$ws.Range("A6:H6").Delete()
$ws.Range("H1:H5").Delete()
